$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'283.49"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").Value = "'20.81"
$ws.Range("D3").ClearFormats()
$ws.Range("D4").Value = "'6.223"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").Value = "'0.06195"
$ws.Range("D5").ClearFormats()
$ws.Range("D7").Value = "'6.563"
$ws.Range("D7").ClearFormats()
$ws.Range("D10").Value = "'0.01390"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").Value = "'0.1652"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "'0.08331"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").Value = "'0.03624"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").Value = "'0.03133"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").Value = "'0.09137"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").Value = "'3.698"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").Value = "'0.001638"
$ws.Range("D17").ClearFormats()
$ws.Range("D19").Value = "'0.006477"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").Value = "'0.006198"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").Value = "'0.001067"
$ws.Range("D21").ClearFormats()
$ws.Range("D23").Value = "'3.820"
$ws.Range("D23").ClearFormats()
$ws.Range("D24").Value = "'2.322"
$ws.Range("D24").ClearFormats()
$ws.Range("D40").Value = "'0.04709"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").Value = "'0.007033"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").Value = "'0.1104"
$ws.Range("D42").ClearFormats()
$ws.Range("D44").Value = "'0.01130"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").Value = "'0.00006364"
$ws.Range("D45").ClearFormats()
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").ClearFormats()
$ws.Range("D47").Value = "'0.9994"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").Value = "'0.002728"
$ws.Range("D48").ClearFormats()
